# Aula 47 - Excluindo cargos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 45
$dstRow = 46

# Copy the formatting (styles) of the previous data row so the new row
# keeps the same look (same font/fill/wrap as the rest of the table).
$ws.Range("B$srcRow`:E$srcRow").Copy()
$ws.Range("B$dstRow`:E$dstRow").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($dstRow, 2).Value = 47
$ws.Cells.Item($dstRow, 3).Value = "9. Cargo: Controller & View"
$ws.Cells.Item($dstRow, 4).Value = "47. Excluindo cargos"
$ws.Cells.Item($dstRow, 5).Value = "2:39`nnova forma de concatenar valores no documento HTML utilizando recursos do próprio thymeleaf`nth:id=""`${#strings.concat('btn_cargos/excluir/', nomeVariavelQlqr.id)}"

$ws.Rows.Item($dstRow).RowHeight = 60

$ws.Application.Goto($ws.Range("D$dstRow"))
